$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.7484373362386563
$ws.Range("E2").Value = 0.7484373362386563

$ws.Range("D3").Value = 0.8487155006798048
$ws.Range("E3").Value = 0.8487155006798048

$ws.Range("D4").Value = [double]"9.565706026693842E-25"
$ws.Range("E4").Value = [double]"9.565706026693842E-25"

$ws.Range("D5").Value = 0.8292039860209397
$ws.Range("E5").Value = 0.8292039860209397

$ws.Range("D6").Value = 0.9962798141428474
$ws.Range("E6").Value = 0.9962798141428474

$ws.Range("D7").Value = [double]"1.011255406600327E-08"
$ws.Range("E7").Value = 0.999999989887446

$ws.Range("D8").Value = 0.9998300371259385
$ws.Range("E8").Value = 0.0001699628740614578

$ws.Range("D9").Value = 0.9607852157971946
$ws.Range("E9").Value = 0.03921478420280544

$ws.Range("D10").Value = 0.9990757685891173
$ws.Range("E10").Value = 0.0009242314108827143

$ws.Range("D11").Value = [double]"1.054564921471774E-07"
$ws.Range("E11").Value = 0.9999998945435079
$ws.Range("F11").Value = 4.514548301696777
